$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2024-10-10 Thursday", $false, $false, $false, $false, $false, $true, 1, $false, "2024-10-11 Friday", 2)

$t = $d.Tables.Item(1)

# Map of (row, col) -> new text for the 5 data rows (1, 5, 9, 13, 17) and 5 columns
$cellValues = @(
    @{Row=1;  Col=1; Text="148÷9=16, 4"},
    @{Row=1;  Col=2; Text="418÷3=139, 1"},
    @{Row=1;  Col=3; Text="129÷8=16, 1"},
    @{Row=1;  Col=4; Text="527÷2=263, 1"},
    @{Row=1;  Col=5; Text="274÷2=137, 0"},

    @{Row=5;  Col=1; Text="678÷8=84, 6"},
    @{Row=5;  Col=2; Text="998÷4=249, 2"},
    @{Row=5;  Col=3; Text="343÷6=57, 1"},
    @{Row=5;  Col=4; Text="581÷8=72, 5"},
    @{Row=5;  Col=5; Text="959÷5=191, 4"},

    @{Row=9;  Col=1; Text="350÷6=58, 2"},
    @{Row=9;  Col=2; Text="558÷7=79, 5"},
    @{Row=9;  Col=3; Text="933÷8=116, 5"},
    @{Row=9;  Col=4; Text="188÷2=94, 0"},
    @{Row=9;  Col=5; Text="302÷4=75, 2"},

    @{Row=13; Col=1; Text="210÷8=26, 2"},
    @{Row=13; Col=2; Text="636÷7=90, 6"},
    @{Row=13; Col=3; Text="907÷6=151, 1"},
    @{Row=13; Col=4; Text="787÷3=262, 1"},
    @{Row=13; Col=5; Text="412÷8=51, 4"},

    @{Row=17; Col=1; Text="791÷5=158, 1"},
    @{Row=17; Col=2; Text="600÷7=85, 5"},
    @{Row=17; Col=3; Text="868÷3=289, 1"},
    @{Row=17; Col=4; Text="608÷7=86, 6"},
    @{Row=17; Col=5; Text="538÷5=107, 3"}
)

foreach ($cv in $cellValues) {
    $cell = $t.Cell($cv.Row, $cv.Col)
    $r = $cell.Range
    $r.SetRange($r.Start, $r.End - 1)
    $r.Text = $cv.Text
}
